$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Range("B15").Value = 96526
$ws.Range("D15").Value = 99.9
$ws.Range("E15").Value = ""

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8"
$ws.Range("B16").Value = 328411
$ws.Range("D16").Value = 99.9
$ws.Range("E16").Value = ""

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6"
$ws.Range("B17").Value = 143808
$ws.Range("D17").Value = 99.9
$ws.Range("E17").Value = ""

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4"
$ws.Range("B18").Value = 287148
$ws.Range("D18").Value = 99.9
$ws.Range("E18").Value = ""

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Range("B19").Value = 69578
$ws.Range("D19").Value = 99.9
$ws.Range("E19").Value = ""

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B20").Value = 10661
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = ""

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B21").Value = 56018
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = ""

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B22").Value = 34244
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = ""

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11"
$ws.Range("B23").Value = 67111
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = ""

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1"
$ws.Range("B24").Value = 13016
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = ""

$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B25").Value = 442178
$ws.Range("D25").Value = 99.9
$ws.Range("E25").Value = "2024-11-10"

$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Range("B26").Value = 18721
$ws.Range("D26").Value = 99.9
$ws.Range("E26").Value = "2024-07-23"

$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7"
$ws.Range("B27").Value = 66577
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = "2024-05-09"

$ws.Range("A28").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B28").Value = 14239
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = "2022-05-23"

$ws.Range("A29").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B29").Value = 265400
$ws.Range("D29").Value = 99.9
$ws.Range("E29").Value = "2022-05-01"

$ws.Range("A30").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B30").Value = 77849
$ws.Range("D30").Value = 99.9
$ws.Range("E30").Value = "2021-08-18"

$ws.Range("A31").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1"
$ws.Range("B31").Value = 26241
$ws.Range("D31").Value = 100
$ws.Range("E31").Value = "2021-01-19"

$ws.Range("A32").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1"
$ws.Range("B32").Value = 15730
$ws.Range("D32").Value = 99.9
$ws.Range("E32").Value = "2020-09-28"

$ws.Range("A33").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B33").Value = 59673
$ws.Range("D33").Value = 100
$ws.Range("E33").Value = "2020-08-05"

$ws.Range("A34").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B34").Value = 113652
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = "2019-12-14"

$ws.Range("A35").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2"
$ws.Range("B35").Value = 88435
$ws.Range("D35").Value = 99.9
$ws.Range("E35").Value = "2019-08-31"

$ws.Range("A36").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2"
$ws.Range("B36").Value = 46270
$ws.Range("D36").Value = 100
$ws.Range("E36").Value = "2019-04-23"
